$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 783.35486
$ws.Range("I15").Value = 783.35486
$ws.Range("K15").Value = 2350.06458
$ws.Range("M15").Value = -2181.06458

$ws.Range("H70").Value = 48348.188
$ws.Range("I70").Value = 2550
$ws.Range("K70").Value = 7650
$ws.Range("M70").Value = -7380

$ws.Range("H73").Value = 48348.188
$ws.Range("I73").Value = 2550
$ws.Range("K73").Value = 7650
$ws.Range("M73").Value = -6714

$ws.Range("H86").Value = 2458.6
$ws.Range("I86").Value = 1455.2858
$ws.Range("K86").Value = 1455.2858
$ws.Range("M86").Value = -332.2858000000001

$ws.Range("H89").Value = 2458.6
$ws.Range("I89").Value = 1455.2858
$ws.Range("K89").Value = 7276.429
$ws.Range("M89").Value = -1660.429

$ws.Range("H98").Value = 5758.778
$ws.Range("I98").Value = 5228.625
$ws.Range("K98").Value = 5228.625
$ws.Range("M98").Value = -3730.625

$ws.Range("H113").Value = 2875
$ws.Range("I113").Value = 2875
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2875
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 379
$ws.Range("N113").ClearContents()

$ws.Range("H116").Value = 4712.857
$ws.Range("I116").Value = 4712.857
$ws.Range("K116").Value = 4712.857
$ws.Range("M116").Value = -1270.857

$ws.Range("H122").Value = 5758.778
$ws.Range("I122").Value = 5228.625
$ws.Range("K122").Value = 15685.875
$ws.Range("M122").Value = -13235.875

$ws.Range("H138").Value = 11775.074
$ws.Range("I138").Value = 7995
$ws.Range("K138").Value = 23985
$ws.Range("M138").Value = -18845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 2000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1623
$ws.Range("N45").ClearContents()

$ws.Range("H122").Value = 478791.62
$ws.Range("I122").Value = 835815.25
$ws.Range("K122").Value = 2507445.75
$ws.Range("M122").Value = -2504995.75

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 225.33333
$ws.Range("I22").Value = 225.33333
$ws.Range("K22").Value = 225.33333
$ws.Range("M22").Value = 124.66667

$ws.Range("H31").Value = 2663.4707
$ws.Range("I31").Value = 3194.1667
$ws.Range("J31").Value = 2374
$ws.Range("K31").Value = 3194.1667
$ws.Range("L31").Value = 2374
$ws.Range("M31").Value = -2899.1667
$ws.Range("N31").Value = -2964

$ws.Range("H34").Value = 2663.4707
$ws.Range("I34").Value = 3194.1667
$ws.Range("J34").Value = 2374
$ws.Range("K34").Value = 3194.1667
$ws.Range("L34").Value = 2374
$ws.Range("M34").Value = -2992.1667
$ws.Range("N34").Value = -2778

$ws.Range("H58").Value = 2735.6538
$ws.Range("J58").Value = 3686.2144
$ws.Range("L58").Value = 3686.2144
$ws.Range("N58").Value = -4092.2144

$ws.Range("H62").Value = 70163.836
$ws.Range("J62").Value = 103247
$ws.Range("L62").Value = 103247
$ws.Range("N62").Value = -104495

$ws.Range("H65").Value = 70163.836
$ws.Range("J65").Value = 103247
$ws.Range("L65").Value = 516235
$ws.Range("N65").Value = -522475

$ws.Range("H132").Value = 4109.1055
$ws.Range("I132").Value = 2672.1333
$ws.Range("K132").Value = 8016.3999
$ws.Range("M132").Value = -5486.3999

$ws.Range("H136").Value = 2735.6538
$ws.Range("J136").Value = 3686.2144
$ws.Range("L136").Value = 11058.6432
$ws.Range("N136").Value = -16158.6432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 334
$ws.Range("I12").Value = 329.64706
$ws.Range("J12").Value = 337.36365
$ws.Range("K12").Value = 988.94118
$ws.Range("L12").Value = 1012.09095
$ws.Range("M12").Value = -815.94118
$ws.Range("N12").Value = -1358.09095

$ws.Range("H38").Value = 95.3125
$ws.Range("I38").Value = 58.7
$ws.Range("J38").Value = 156.33333
$ws.Range("K38").Value = 176.1
$ws.Range("L38").Value = 468.99999
$ws.Range("M38").Value = 170.9
$ws.Range("N38").Value = -1162.99999

$ws.Range("H70").Value = 663.6667
$ws.Range("I70").Value = 663.6667
$ws.Range("K70").Value = 1991.0001
$ws.Range("M70").Value = -1676.0001

$ws.Range("H73").Value = 663.6667
$ws.Range("I73").Value = 663.6667
$ws.Range("K73").Value = 1991.0001
$ws.Range("M73").Value = -899.0001

$ws.Range("H75").Value = 76.75
$ws.Range("I75").Value = 43
$ws.Range("J75").Value = 88
$ws.Range("K75").Value = 129
$ws.Range("L75").Value = 264
$ws.Range("M75").Value = 869
$ws.Range("N75").Value = -2260

$ws.Range("H78").Value = 76.75
$ws.Range("I78").Value = 43
$ws.Range("J78").Value = 88
$ws.Range("K78").Value = 387
$ws.Range("L78").Value = 792
$ws.Range("M78").Value = 4605
$ws.Range("N78").Value = -10776

$ws.Range("H137").Value = 5875.9
$ws.Range("I137").Value = 4952.6
$ws.Range("J137").Value = 6799.2
$ws.Range("K137").Value = 14857.8
$ws.Range("L137").Value = 20397.6
$ws.Range("M137").Value = -9757.800000000001
$ws.Range("N137").Value = -30597.6

$ws.Range("H139").Value = 2095.3333
$ws.Range("I139").Value = 2095.3333
$ws.Range("K139").Value = 6285.999899999999
$ws.Range("M139").Value = -1145.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H122").Value = 32327.334
$ws.Range("I122").Value = 1944.7826
$ws.Range("K122").Value = 5834.3478
$ws.Range("M122").Value = -3384.3478

$ws.Range("H123").Value = 18599.334
$ws.Range("J123").Value = 18599.334
$ws.Range("L123").Value = 18599.334
$ws.Range("N123").Value = -23499.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -887

$ws.Range("H16").Value = 15049.4
$ws.Range("J16").Value = 14082.667
$ws.Range("L16").Value = 14082.667
$ws.Range("N16").Value = -14422.667

$ws.Range("H28").Value = 1000
$ws.Range("I28").Value = 1000
$ws.Range("K28").Value = 1000
$ws.Range("M28").Value = -768

$ws.Range("H37").Value = 1000
$ws.Range("I37").Value = 1000
$ws.Range("K37").Value = 1000
$ws.Range("M37").Value = -893

$ws.Range("H46").Value = 2779.913
$ws.Range("I46").Value = 2077.7334
$ws.Range("J46").Value = 4096.5
$ws.Range("K46").Value = 2077.7334
$ws.Range("L46").Value = 4096.5
$ws.Range("M46").Value = -1889.7334
$ws.Range("N46").Value = -4472.5

$ws.Range("H122").Value = 4582.636
$ws.Range("I122").Value = 4489.8887
$ws.Range("K122").Value = 13469.6661
$ws.Range("M122").Value = -11019.6661

$ws.Range("H136").Value = 2687.75
$ws.Range("I136").Value = 2220.6
$ws.Range("K136").Value = 6661.799999999999
$ws.Range("M136").Value = -4111.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1898.1666
$ws.Range("J122").Value = 10000
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -34900

$ws.Range("H132").Value = 106379.89
$ws.Range("I132").Value = 189964.8
$ws.Range("K132").Value = 569894.3999999999
$ws.Range("M132").Value = -567364.3999999999
Write-Host "done"
